# Add a new "september-2025" worksheet at the end of the workbook, mirroring
# the layout of the existing monthly sheets (a single A1 cell containing the
# tax-revenue summary text for that month).

$wb = $excel.ActiveWorkbook

# Find the current last sheet so the new one can be inserted right after it.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "september-2025"

$newSheet.Range("A1").Value = ": tax revenue                                              71,294            68,173            3,121            4.6%"
